$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Business Exceptions")

$textValue = "text"
$longValue = "CHANGES - SOX...eml on (selected date) does not contain Server Name "

$headerRows = @(71, 90, 109)

for ($r = 71; $r -le 124; $r++) {
    if ($headerRows -contains $r) {
        $ws.Cells.Item($r, 1).Value = $textValue
        $ws.Cells.Item($r, 2).Value = 1
    } else {
        $ws.Cells.Item($r, 1).Value = $longValue
    }
}
